# Weekly fruit/vegetable price update.
# Insert two new report rows (most recent week) above the existing row 11
# data block, pushing the previous rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 11 - this shifts rows 11:18 down to 13:20
# (including their formatting, e.g. the date style on column D).
$ws.Range("A11:R12").Insert()

# Row 11: new "Cultivar XV región" / "Primera" entry
$ws.Cells.Item(11, 1).Value = 1
$ws.Cells.Item(11, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(11, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(11, 4).Value = 44533
$ws.Cells.Item(11, 5).Value = 15
$ws.Cells.Item(11, 6).Value = 100112043
$ws.Cells.Item(11, 7).Value = "Pepino dulce"
$ws.Cells.Item(11, 8).Value = "Cultivar XV región"
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 100
$ws.Cells.Item(11, 11).Value = 6000
$ws.Cells.Item(11, 12).Value = 7000
$ws.Cells.Item(11, 13).Value = 6500
$ws.Cells.Item(11, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(11, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(11, 16).Value = 650
$ws.Cells.Item(11, 17).Value = 10
$ws.Cells.Item(11, 18).Value = "Hortaliza"

# Row 12: new "Cultivar XV región" / "Segunda" entry
$ws.Cells.Item(12, 1).Value = 1
$ws.Cells.Item(12, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(12, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(12, 4).Value = 44533
$ws.Cells.Item(12, 5).Value = 15
$ws.Cells.Item(12, 6).Value = 100112043
$ws.Cells.Item(12, 7).Value = "Pepino dulce"
$ws.Cells.Item(12, 8).Value = "Cultivar XV región"
$ws.Cells.Item(12, 9).Value = "Segunda"
$ws.Cells.Item(12, 10).Value = 120
$ws.Cells.Item(12, 11).Value = 4000
$ws.Cells.Item(12, 12).Value = 5000
$ws.Cells.Item(12, 13).Value = 4500
$ws.Cells.Item(12, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(12, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(12, 16).Value = 450
$ws.Cells.Item(12, 17).Value = 10
$ws.Cells.Item(12, 18).Value = "Hortaliza"
